# Generate Report for Handback
# Row 7 ("ee91ac2e-cd55-46ea-91af-fcbe2744574f.md") on both the zh-cn and
# de-de sheets now has a handback: a (stale) target/handback file, a
# handback datetime, and an error explaining the handback version is not
# the latest one.

$wb = $excel.ActiveWorkbook

$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/44c23a7465a5bd2009e3bf7d255b9e95ae914e5e/e2e/ee91ac2e-cd55-46ea-91af-fcbe2744574f.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/03970fce9aa911d1197b98ded044f6a05aa9afdd/e2e/ee91ac2e-cd55-46ea-91af-fcbe2744574f.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/44c23a7465a5bd2009e3bf7d255b9e95ae914e5e/e2e/ee91ac2e-cd55-46ea-91af-fcbe2744574f.md."

function Update-HandbackRow($SheetName, $HandbackDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Latest Target File (I7): now resolved to the source .md, same
    # display/link pattern as the other rows' "Latest Target File" cell.
    $ws.Range("I7").Value = "ee91ac2e-cd55-46ea-91af-fcbe2744574f.md"
    $ws.Range("I7").Style = "HyperLink"
    $ws.Hyperlinks.Add($ws.Range("I7"), $handbackUrl, "", "", "ee91ac2e-cd55-46ea-91af-fcbe2744574f.md")

    # Latest Handback File (J7): same xlf file referenced by Latest Handoff
    # File (G7) for this row/locale.
    $ws.Range("J7").Value = $ws.Range("G7").Value2

    # Latest Handback DateTime (K7)
    $ws.Range("K7").Value = $HandbackDateTime

    # Error Detail (P7)
    $ws.Range("P7").Value = $errorDetail
}

Update-HandbackRow "zh-cn" "2016-09-01 07:04:38"
Update-HandbackRow "de-de" "2016-09-01 07:04:45"
